$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 2.98
$ws.Range("I2").Value = 3.1
$ws.Range("L2").Value = 1.57
$ws.Range("Y2").Value = 1000
$ws.Range("S3").Value = 1.73
$ws.Range("U3").Value = 1.98
$ws.Range("AF3").Value = 12.5
$ws.Range("AG3").Value = 16.5
$ws.Range("AJ3").Value = 12.5
$ws.Range("AK3").Value = 16.5
$ws.Range("G4").Value = 1.97
$ws.Range("M4").Value = 1.12
$ws.Range("P4").Value = 1.52
$ws.Range("W4").Value = 2.02
$ws.Range("G5").Value = 2.6
$ws.Range("I5").Value = 3.45
$ws.Range("J5").Value = 3.15
$ws.Range("N5").Value = 3
$ws.Range("O5").Value = 1.47
$ws.Range("S5").Value = 4.7
$ws.Range("U5").Value = 1.93
$ws.Range("W5").Value = 1.62
$ws.Range("AK5").Value = 34
$ws.Range("AM5").Value = 150
$ws.Range("AN5").Value = 34
$ws.Range("G6").Value = 1.63
$ws.Range("K6").Value = 5.6
$ws.Range("R6").Value = 1.93
$ws.Range("S6").Value = 1.91
$ws.Range("X6").Value = 42
$ws.Range("AC6").Value = 14.5
$ws.Range("AL6").Value = 21
$ws.Range("H7").Value = 1.76
$ws.Range("I7").Value = 1.86
$ws.Range("J7").Value = 3.75
$ws.Range("K7").Value = 4.1
$ws.Range("L7").Value = 1.4
$ws.Range("N7").Value = 3.55
$ws.Range("P7").Value = 1.9
$ws.Range("Q7").Value = 1.96
$ws.Range("S7").Value = 3.35
$ws.Range("U7").Value = 1.97
$ws.Range("V7").Value = 2.16
$ws.Range("Z7").Value = 11.5
$ws.Range("AA7").Value = 980
$ws.Range("AC7").Value = 9.199999999999999
$ws.Range("AE7").Value = 980
$ws.Range("AF7").Value = 46
$ws.Range("AK7").Value = 75
$ws.Range("AO7").Value = 13.5
$ws.Range("S8").Value = 4.2
$ws.Range("I9").Value = 10
$ws.Range("K9").Value = 4.7
$ws.Range("N9").Value = 3.55
$ws.Range("O9").Value = 1.37
$ws.Range("P9").Value = 1.87
$ws.Range("Q9").Value = 2.1
$ws.Range("S9").Value = 3.95
$ws.Range("F10").Value = 1.59
$ws.Range("G10").Value = 1.6
$ws.Range("L10").Value = 1.37
$ws.Range("R10").Value = 1.45
$ws.Range("T10").Value = 1.9
$ws.Range("U10").Value = 2.06
$ws.Range("V10").Value = 1.18
$ws.Range("W10").Value = 2.66
$ws.Range("X10").Value = 18
$ws.Range("Z10").Value = 55
$ws.Range("AA10").Value = 180
$ws.Range("AB10").Value = 8.800000000000001
$ws.Range("AC10").Value = 10
$ws.Range("AE10").Value = 85
$ws.Range("AF10").Value = 9.199999999999999
$ws.Range("L11").Value = 1.48
$ws.Range("X11").Value = 13
